# FixArt_candidates.xlsx — "fixed figs in Ch07; avoid deprecated effects:: warning in Ch11"
#
# Ch07 rows (Fig 7.6 "arth-cond1", Fig 7.7 "arth-cond2", Fig 7.15 "donner1-cond1",
# Fig 7.16 "donner1-cond3") previously had the Resolution column (D) filled in with
# the TODO note "FIXME: should use theme_bw()" (shown with the red/"Bad" highlight
# style). Those plots have now actually been fixed, so the note is replaced with
# "fixed + theme_bw()" and re-styled with Excel's built-in green "Good" cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$fixedCells = "D39", "D40", "D43", "D44"

foreach ($cellRef in $fixedCells) {
    $ws.Range($cellRef).Value = "fixed + theme_bw()"
}

foreach ($cellRef in $fixedCells) {
    $ws.Range($cellRef).Style = "Good"
}

# Scroll the sheet down and leave the selection on the last-fixed cell, as it was
# left after making the edit.
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D44").Select()
